$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.959.97"
$ws.Range("E2").Value = "  -0.74%  "

$ws.Range("D3").Value = "2.450.96"
$ws.Range("E3").Value = "  +0.79%  "

$ws.Range("E4").Value = "  -0.15%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.54"
$ws.Range("E5").Value = "  +0.96%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.76"
$ws.Range("E6").Value = "  -0.81%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.530"
$ws.Range("E8").Value = "  +0.37%  "

$ws.Range("D9").Value = "2.446.90"
$ws.Range("E9").Value = "  +0.81%  "

$ws.Range("E10").Value = "  +2.27%  "

$ws.Range("E11").Value = "  +2.83%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.17"
$ws.Range("E12").Value = "  -0.79%  "

$ws.Range("E13").Value = "  -2.28%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.94"
$ws.Range("E14").Value = "  -1.21%  "

$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000172"
$ws.Range("E15").Value = "  -0.24%  "

$ws.Range("B16").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C16").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D16").Value = "2.887.60"
$ws.Range("E16").Value = "  -0.19%  "

$ws.Range("D17").Value = "61.895.61"
$ws.Range("E17").Value = "  -0.78%  "

$ws.Range("D18").Value = "2.440.42"
$ws.Range("E18").Value = "  +0.15%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.62"
$ws.Range("E19").Value = "  -3.39%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.22"
$ws.Range("E20").Value = "  +1.76%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "325.09"
$ws.Range("E21").Value = "  -1.02%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.07"
$ws.Range("E22").Value = "  -1.04%  "

$ws.Range("E23").Value = "  +0.04%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.92"
$ws.Range("E24").Value = "  -2.12%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.07"
$ws.Range("E25").Value = "  -0.73%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.16"
$ws.Range("E26").Value = "  +3.72%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "586.21"
$ws.Range("E27").Value = "  -6.80%  "

$ws.Range("D28").Value = "2.562.69"
$ws.Range("E28").Value = "  -0.04%  "

$ws.Range("E29").Value = "  +0.05%  "

$ws.Range("D30").Value = "0.0₃0930"
$ws.Range("E30").Value = "  -2.73%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.88"
$ws.Range("E31").Value = "  -1.65%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.37"
$ws.Range("E32").Value = "  -3.77%  "

$ws.Range("E33").Value = "  -0.27%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.134"
$ws.Range("E34").Value = "  -3.11%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.77"
$ws.Range("E36").Value = "  -3.69%  "

$ws.Range("E37").Value = "  -0.54%  "

$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.40"
$ws.Range("E38").Value = "  -3.67%  "

$ws.Range("B39").Value = "Monero"
$ws.Range("C39").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "151.86"
$ws.Range("E39").Value = "  +3.10%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.32"
$ws.Range("E40").Value = "  -0.44%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.16"
$ws.Range("E41").Value = "  -1.61%  "

$ws.Range("B42").Value = "OKB"
$ws.Range("C42").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "42.28"
$ws.Range("E42").Value = "  -0.29%  "

$ws.Range("B43").Value = "USDe"
$ws.Range("C43").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.999"
$ws.Range("E43").Value = "  +0.02%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.68"
$ws.Range("E44").Value = "  -3.13%  "

$ws.Range("E45").Value = "  -3.78%  "

$ws.Range("D46").Value = "0.0₆0277"
$ws.Range("E46").Value = "  +19.65%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "141.50"
$ws.Range("E47").Value = "  -1.60%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.57"
$ws.Range("E48").Value = "  -2.78%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.598"
$ws.Range("E49").Value = "  +0.42%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0510"
$ws.Range("E50").Value = "  -2.00%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.62"
$ws.Range("E51").Value = "  +0.40%  "
